# Commit: "change status to published"
#
# For every glossary row that was previously curated/"Discussed", mark the
# curator as having reviewed it (append/set "JH" in the Curator column) and
# flip the Curation status from "Discussed" to "Published". These are all
# the rows that share the moccasin-style background (style index 2 /
# fillId 2); that shared fill is also recolored from moccasin (#FFE4B5) to
# aquamarine (#7FFFD4) to visually reflect the new "Published" status.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Curation status is "Discussed" and should become "Published".
$rows = @(2,3,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,34,35,36,37,40,41,42,43,44)

foreach ($r in $rows) {
    $gCell = $ws.Cells.Item($r, 7)   # column G = Curator
    $hCell = $ws.Cells.Item($r, 8)   # column H = Curation status

    $curCurator = $gCell.Text
    if ($curCurator -and $curCurator.Trim().Length -gt 0) {
        $gCell.Value = "$curCurator; JH"
    } else {
        $gCell.Value = "JH"
    }

    $hCell.Value = "Published"
}

# Recolor the shared background fill used by these rows from moccasin
# (00FFE4B5) to aquamarine (007FFFD4). The rows form four contiguous
# blocks (A2:J3, A5:J32, A34:J37, A40:J44); set each block's own Interior
# color on its Areas so the engine reuses a single shared style.
$fillRange = $ws.Range("A2:J3,A5:J32,A34:J37,A40:J44")
foreach ($area in $fillRange.Areas) {
    $area.Interior.Color = 13959039   # RGB(127,255,212) = 0x7FFFD4
}

Write-Host "Updated" $rows.Count "rows to Published and recolored their fill."
